# Apply the Jun 8 2024 cryptos-list refresh (prices + 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get numeric-looking text (e.g. "0.999", "684.31").
# Excel auto-converts a numeric-looking string typed into .Value into a
# real number, which would silently turn "1.00" into 1. Forcing the
# cell to Text format first keeps them as literal strings, matching the
# source data (which stores every Price/Volume cell as text).
$textPriceCells = @(
    "D4", "D5", "D6", "D9", "D10", "D11", "D12", "D13",
    "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D26",
    "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35",
    "D36", "D39", "D42", "D45", "D46", "D47", "D48", "D49",
    "D51"
)
foreach ($ref in $textPriceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Cell value updates ---
$ws.Range("D2").Value = "69.356.56"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "3.683.88"
$ws.Range("E3").Value = "  -3.92%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "684.31"
$ws.Range("D6").Value = "160.30"
$ws.Range("E6").Value = "  -7.05%  "
$ws.Range("D7").Value = "3.683.60"
$ws.Range("E7").Value = "  -4.03%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -6.37%  "
$ws.Range("D10").Value = "0.146"
$ws.Range("E10").Value = "  -9.56%  "
$ws.Range("D11").Value = "7.28"
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  -10.26%  "
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").Value = "  -7.86%  "
$ws.Range("D14").Value = "4.301.48"
$ws.Range("E14").Value = "  -4.01%  "
$ws.Range("D15").Value = "32.67"
$ws.Range("E15").Value = "  -10.74%  "
$ws.Range("D16").Value = "3.683.58"
$ws.Range("E16").Value = "  -3.63%  "
$ws.Range("D17").Value = "69.367.30"
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "15.89"
$ws.Range("E19").Value = "  -9.88%  "
$ws.Range("D20").Value = "6.46"
$ws.Range("E20").Value = "  -11.31%  "
$ws.Range("D21").Value = "470.49"
$ws.Range("E21").Value = "  -9.13%  "
$ws.Range("D22").Value = "9.90"
$ws.Range("E22").Value = "  -6.50%  "
$ws.Range("D23").Value = "0.652"
$ws.Range("E23").Value = "  -9.54%  "
$ws.Range("D24").Value = "79.54"
$ws.Range("E24").Value = "  -5.63%  "
$ws.Range("D25").Value = "3.825.38"
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "0.0000127"
$ws.Range("E26").Value = "  -11.78%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "11.10"
$ws.Range("E28").Value = "  -13.19%  "
$ws.Range("D29").Value = "9.11"
$ws.Range("E29").Value = "  -12.61%  "
$ws.Range("D30").Value = "2.69"
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").Value = "  -14.08%  "
$ws.Range("D32").Value = "6.72"
$ws.Range("E32").Value = "  -9.57%  "
$ws.Range("D33").Value = "2.03"
$ws.Range("E33").Value = "  -10.37%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "26.72"
$ws.Range("E35").Value = "  -9.18%  "
$ws.Range("D36").Value = "0.162"
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("D37").Value = "3.647.67"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("E38").Value = "  -11.70%  "
$ws.Range("D39").Value = "6.11"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("E40").Value = "  -8.93%  "
$ws.Range("D42").Value = "0.0908"
$ws.Range("E42").Value = "  -10.66%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  -7.56%  "
$ws.Range("D45").Value = "164.20"
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("D46").Value = "48.08"
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("D47").Value = "29.40"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.71"
$ws.Range("E48").Value = "  -18.04%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("E50").Value = "  -12.45%  "
$ws.Range("D51").Value = "1.10"
$ws.Range("E51").Value = "  -5.17%  "

# Restore default (unstyled) formatting now that the text is committed,
# so these cells end up styleless just like in the source workbook.
foreach ($ref in $textPriceCells) {
    $ws.Range($ref).Style = "Normal"
}
